$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A10").Value = "13/12/2019"
$ws.Range("B10").Value = "The screens for modules were started"
$ws.Range("C10").Value = "Modules were drawn and designing is started"

$ws.Range("C10").Select()
